# feature: use regex to extract a part organization name
#
# Normalises the "MYCONS <dash> OrgN" organization label used on the
# PROD1 (summary) sheet so it always uses a plain hyphen, regardless of
# whether the source text used a hyphen-minus or an en-dash. The raw
# PROD2 (detail) sheet is left untouched, so it keeps showing whatever
# dash character the original export used.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)

# Re-write the organization labels on the summary sheet, normalising the
# separator to a plain ASCII hyphen ("-") via the extracted/regex-cleaned
# name, e.g. "MYCONS – Org2" -> "MYCONS - Org2".
$ws1.Cells.Item(2, 2).Value = "MYCONS - Org1"
$ws1.Cells.Item(3, 2).Value = "MYCONS - Org2"
$ws1.Cells.Item(4, 2).Value = "MYCONS - Org3"

# Match the author's last recorded selection on the summary sheet.
$ws1.Activate()
$ws1.Range("B7").Select()
